$d = $word.ActiveDocument

# The first (Title) paragraph currently holds the full combined title text:
#   "Let's talk about Thurstone & Co.: An information-theoretical model for
#    comparative judgments, and its statistical translation"
# It needs to be split into:
#   Title paragraph:    "Let's talk about Thurstone & Co. "
#   Subtitle paragraph: "An information-theoretical model for comparative
#                         judgments, and its statistical translation"

$newTitle = "Let’s talk about Thurstone & Co. "
$subtitle = "An information-theoretical model for comparative judgments, and its statistical translation"

$titlePara = $d.Paragraphs.First

# Replace the title paragraph's text with just the shortened title.
$titlePara.Range.Text = $newTitle

# Insert a brand-new paragraph right after the Title paragraph for the
# subtitle text. We build it via raw OOXML so the emitted <w:t> element keeps
# xml:space="preserve" (matching the document's existing convention) even
# though the subtitle text has no leading/trailing whitespace of its own.
$endOfTitle = $titlePara.Range
$endOfTitle.Collapse(0) | Out-Null
$endOfTitle.InsertParagraphAfter()

$subtitlePara = $d.Paragraphs.Item(2)
$subtitleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t xml:space="preserve">' + $subtitle + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$subtitlePara.Range.InsertXML($subtitleXml)
